# Applies the row-content changes described in the diff:
#  - Rows 34 and 35 swap their entire content.
#  - Rows 37, 38, 39 rotate their content: new row37 = old row39,
#    new row38 = old row37, new row39 = old row38.
# The row objects themselves (row number, any row-level formatting) stay
# put; only the cell values move between rows, exactly as in the diff.
#
# NOTE: this sheet stores a handful of numeric-looking values (e.g. the
# "Antal" column I, values like "1"/"2") as TEXT, not numbers. Excel's COM
# layer will happily "helpfully" reinterpret a plain numeric string as a
# real number/date the moment it is written back through .Value, so those
# columns are forced to text format (NumberFormat "@") before the write to
# keep them as text, matching the source data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colCount = $ws.UsedRange.Columns.Count   # number of used columns (A..AY == 51)

# Columns that hold genuine numbers (Id, Taxonsorteringsordning, TaxonId,
# Ost, Nord, Noggrannhet) -- everything else on this sheet is text (or
# boolean), even when it looks numeric (e.g. column I "Antal": "1", "2").
$numericCols = @(1, 2, 5, 17, 18, 19)   # A, B, E, Q, R, S

function Get-RowValues($ws, $row, $colCount) {
    $vals = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $vals += ,$ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

# Writes $val into the cell at ($row, $col), preserving the data's
# original type:
#  - $null                             -> clear the cell
#  - non-numeric column AND the value
#    looks like a bare number          -> force text format first so Excel
#                                          doesn't silently re-interpret a
#                                          numeric-looking value (e.g. the
#                                          "Antal" column stores "1"/"2" as
#                                          TEXT, never as a real number)
#  - everything else                   -> plain assignment
function Set-CellSmart($ws, $row, $col, $val, $numericCols) {
    $cell = $ws.Cells.Item($row, $col)
    if ($null -eq $val) {
        $cell.ClearContents() | Out-Null
        return
    }
    if (-not ($numericCols -contains $col)) {
        $looksNumeric = "$val" -match '^-?\d+(\.\d+)?$'
        if ($looksNumeric) {
            $cell.NumberFormat = "@"
        }
    }
    $cell.Value = $val
}

function Set-RowValues($ws, $row, $colCount, $oldVals, $newVals, $numericCols) {
    for ($c = 1; $c -le $colCount; $c++) {
        $old = $oldVals[$c - 1]
        $new = $newVals[$c - 1]
        if ($old -ne $new) {
            Set-CellSmart $ws $row $c $new $numericCols
        }
    }
}

# --- Snapshot the current ("before") content of every affected row ---
$orig34 = Get-RowValues $ws 34 $colCount
$orig35 = Get-RowValues $ws 35 $colCount

$orig37 = Get-RowValues $ws 37 $colCount
$orig38 = Get-RowValues $ws 38 $colCount
$orig39 = Get-RowValues $ws 39 $colCount

# --- Swap rows 34 and 35 ---
Set-RowValues $ws 34 $colCount $orig34 $orig35 $numericCols
Set-RowValues $ws 35 $colCount $orig35 $orig34 $numericCols

# --- Rotate rows 37 -> 38 -> 39 -> 37 ---
Set-RowValues $ws 37 $colCount $orig37 $orig39 $numericCols
Set-RowValues $ws 38 $colCount $orig38 $orig37 $numericCols
Set-RowValues $ws 39 $colCount $orig39 $orig38 $numericCols
